# Widen the distance between the two columns of the (single) section
# from 851 twips (42.55 pt) to 709 twips (35.45 pt).
#
# Word's object model expresses this via PageSetup.TextColumns.Spacing,
# measured in points (1 pt = 20 twips), which maps directly onto the
# <w:cols w:space="..."/> attribute in word/document.xml.

$d = $word.ActiveDocument

foreach ($section in $d.Sections) {
    $section.PageSetup.TextColumns.Spacing = 709 / 20
}
